$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as Text so values like "235.78" are not
# auto-converted to numbers, matching the original inline-string text cells.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Rows whose Coin/Link swapped position with a neighboring row ---
$ws.Range("B17").Value = "Dai"
$ws.Range("C17").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D17").Value = "0.9997"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "13.20"
$ws.Range("E18").Value = "  -2.02%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.281"
$ws.Range("E23").Value = "  -1.45%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "0.9014"
$ws.Range("E39").Value = "  -2.37%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "2.048"
$ws.Range("E40").Value = "  +0.20%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  +0.61%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4241"
$ws.Range("E44").Value = "  -3.35%  "

# --- Rows with only Price/Volume(1h) updates ---
$ws.Range("D2").Value = "30.218.86"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.855.31"
$ws.Range("E3").Value = "  -2.30%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "235.78"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D7").Value = "0.4781"
$ws.Range("E7").Value = "  -2.71%  "
$ws.Range("D8").Value = "0.2803"
$ws.Range("E8").Value = "  -4.44%  "
$ws.Range("D9").Value = "0.06477"
$ws.Range("E9").Value = "  -3.42%  "
$ws.Range("D10").Value = "1.858.07"
$ws.Range("E10").Value = "  -3.03%  "
$ws.Range("D11").Value = "0.07348"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "16.19"
$ws.Range("E12").Value = "  -4.57%  "
$ws.Range("D13").Value = "5.096"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "87.15"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "0.6463"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "30.165.37"
$ws.Range("E16").Value = "  -0.92%  "
$ws.Range("D19").Value = "0.000007616"
$ws.Range("E19").Value = "  -2.89%  "
$ws.Range("D20").Value = "225.41"
$ws.Range("E20").Value = "  +17.40%  "
$ws.Range("D21").Value = "2.105.82"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("D24").Value = "6.061"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "9.227"
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("D26").Value = "163.88"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "18.41"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("D28").Value = "1.918"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "1.444"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").Value = "0.09193"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "4.244"
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("D32").Value = "3.956"
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("D33").Value = "0.05008"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "0.7349"
$ws.Range("E34").Value = "  -0.99%  "
$ws.Range("D35").Value = "1.147"
$ws.Range("E35").Value = "  +4.20%  "
$ws.Range("D36").Value = "2.691"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").Value = "0.01829"
$ws.Range("E37").Value = "  +0.98%  "
$ws.Range("D38").Value = "2.600"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D41").Value = "5.950"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").Value = "106.24"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("D45").Value = "0.1319"
$ws.Range("E45").Value = "  -3.84%  "
$ws.Range("D46").Value = "7.366"
$ws.Range("E46").Value = "  -3.21%  "
$ws.Range("D47").Value = "1.539"
$ws.Range("E47").Value = "  +9.53%  "
$ws.Range("D48").Value = "64.10"
$ws.Range("E48").Value = "  -7.21%  "
$ws.Range("D49").Value = "8.770"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").Value = "34.08"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").Value = "0.05663"
$ws.Range("E51").Value = "  -2.92%  "
